$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear cells N11 and N12 (remove their values entirely)
$ws.Range("N11").ClearContents()
$ws.Range("N12").ClearContents()

# Update the selected cell/range on the sheet view
$ws.Range("R13").Select()
